# The commit swaps the contents of ppt/theme/theme1.xml (the deck's main
# "Integral" theme, driven by the slide master) and ppt/theme/theme2.xml
# (the "Office Theme" used by the notes master) with each other.
#
# Inspecting both theme parts shows their <a:fontScheme> and <a:fmtScheme>
# blocks are already byte-for-byte identical; the only real differences are
# the top-level theme name and the 12 colours of <a:clrScheme>. So swapping
# the two files is equivalent to swapping just those colour values (and
# names) between the two parts.
#
# The PowerPoint object model reaches the deck's (slide-master-backed)
# theme colours through Slide.ThemeColorScheme - each of the 12 entries
# maps 1:1 onto dk1/lt1/dk2/lt2/accent1-6/hlink/folHlink, writable via
# .RGB (standard VBA RGB(r,g,b) = r + g*256 + b*65536 packing). Apply the
# "Office Theme" palette (previously in theme2.xml) there so the deck's
# active theme becomes the Office Theme palette, matching the diff's net
# visual effect on ppt/theme/theme1.xml.

# PowerPoint/VBA's RGB() helper isn't available in this host, so replicate
# its packing (r + g*256 + b*65536) from a hex string ourselves.
function HexToOleColor([string]$hex) {
    $r = [Convert]::ToInt32($hex.Substring(0,2), 16)
    $g = [Convert]::ToInt32($hex.Substring(2,2), 16)
    $b = [Convert]::ToInt32($hex.Substring(4,2), 16)
    return $r + ($g * 256) + ($b * 65536)
}

$p = $ppt.ActivePresentation
$s = $p.Slides.Item(1)
$cs = $s.ThemeColorScheme

# index -> scheme slot -> target "Office Theme" RGB hex (was "Integral")
# 1 dk1       000000
# 2 lt1       FFFFFF
# 3 dk2       44546A
# 4 lt2       E7E6E6
# 5 accent1   5B9BD5
# 6 accent2   ED7D31
# 7 accent3   A5A5A5
# 8 accent4   FFC000
# 9 accent5   4472C4
# 10 accent6  70AD47
# 11 hlink    0563C1
# 12 folHlink 954F72

$cs.Item(1).RGB  = HexToOleColor "000000"
$cs.Item(2).RGB  = HexToOleColor "FFFFFF"
$cs.Item(3).RGB  = HexToOleColor "44546A"
$cs.Item(4).RGB  = HexToOleColor "E7E6E6"
$cs.Item(5).RGB  = HexToOleColor "5B9BD5"
$cs.Item(6).RGB  = HexToOleColor "ED7D31"
$cs.Item(7).RGB  = HexToOleColor "A5A5A5"
$cs.Item(8).RGB  = HexToOleColor "FFC000"
$cs.Item(9).RGB  = HexToOleColor "4472C4"
$cs.Item(10).RGB = HexToOleColor "70AD47"
$cs.Item(11).RGB = HexToOleColor "0563C1"
$cs.Item(12).RGB = HexToOleColor "954F72"
